$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1866
$endRow = 1934
$numRows = $endRow - $startRow + 1
$data = New-Object 'object[,]' $numRows,6

$data[0,0] = "2025-06-09 09:01:52"
$data[0,1] = 516
$data[0,2] = 18.08
$data[0,3] = 789.9400000000001
$data[0,4] = 23.9
$data[0,5] = 73.09999999999999
$data[1,0] = "2025-06-09 09:01:53"
$data[1,1] = 516
$data[1,2] = 18.08
$data[1,3] = 789.8200000000001
$data[1,4] = 23.9
$data[1,5] = 73.09999999999999
$data[2,0] = "2025-06-09 09:01:54"
$data[2,1] = 516
$data[2,2] = 18.57
$data[2,3] = 789.75
$data[2,4] = 23.9
$data[2,5] = 73.2
$data[3,0] = "2025-06-09 09:01:56"
$data[3,1] = 515
$data[3,2] = 18.08
$data[3,3] = 16.66
$data[3,4] = 23.9
$data[3,5] = 73.2
$data[4,0] = "2025-06-09 09:01:57"
$data[4,1] = 516
$data[4,2] = 18.08
$data[4,3] = 28.76
$data[4,4] = 23.9
$data[4,5] = 73.2
$data[5,0] = "2025-06-09 09:01:58"
$data[5,1] = 516
$data[5,2] = 18.08
$data[5,3] = 18.24
$data[5,4] = 23.9
$data[5,5] = 73.2
$data[6,0] = "2025-06-09 09:01:59"
$data[6,1] = 515
$data[6,2] = 18.08
$data[6,3] = 789.84
$data[6,4] = 23.9
$data[6,5] = 73.3
$data[7,0] = "2025-06-09 09:02:00"
$data[7,1] = 516
$data[7,2] = 18.08
$data[7,3] = 789.99
$data[7,4] = 23.9
$data[7,5] = 73.3
$data[8,0] = "2025-06-09 09:02:02"
$data[8,1] = 516
$data[8,2] = 18.08
$data[8,3] = 37.23
$data[8,4] = 23.9
$data[8,5] = 73.3
$data[9,0] = "2025-06-09 09:02:03"
$data[9,1] = 516
$data[9,2] = 18.08
$data[9,3] = 25.93
$data[9,4] = 23.9
$data[9,5] = 73.3
$data[10,0] = "2025-06-09 09:02:04"
$data[10,1] = 515
$data[10,2] = 18.08
$data[10,3] = 789.85
$data[10,4] = 23.9
$data[10,5] = 73.40000000000001
$data[11,0] = "2025-06-09 09:02:05"
$data[11,1] = 515
$data[11,2] = 18.08
$data[11,3] = 789.92
$data[11,4] = 23.9
$data[11,5] = 73.40000000000001
$data[12,0] = "2025-06-09 09:02:07"
$data[12,1] = 516
$data[12,2] = 18.08
$data[12,3] = 25.36
$data[12,4] = 23.9
$data[12,5] = 73.3
$data[13,0] = "2025-06-09 09:02:08"
$data[13,1] = 516
$data[13,2] = 18.08
$data[13,3] = 20.43
$data[13,4] = 23.9
$data[13,5] = 73.3
$data[14,0] = "2025-06-09 09:02:09"
$data[14,1] = 516
$data[14,2] = 18.08
$data[14,3] = 789.85
$data[14,4] = 23.9
$data[14,5] = 73.5
$data[15,0] = "2025-06-09 09:02:11"
$data[15,1] = 516
$data[15,2] = 18.08
$data[15,3] = 789.8200000000001
$data[15,4] = 23.9
$data[15,5] = 73.5
$data[16,0] = "2025-06-09 09:02:12"
$data[16,1] = 516
$data[16,2] = 18.08
$data[16,3] = 28.22
$data[16,4] = 23.8
$data[16,5] = 73.59999999999999
$data[17,0] = "2025-06-09 09:02:13"
$data[17,1] = 516
$data[17,2] = 18.08
$data[17,3] = 23.75
$data[17,4] = 23.8
$data[17,5] = 73.59999999999999
$data[18,0] = "2025-06-09 09:02:15"
$data[18,1] = 516
$data[18,2] = 18.08
$data[18,3] = 24.5
$data[18,4] = 23.8
$data[18,5] = 73.59999999999999
$data[19,0] = "2025-06-09 09:02:16"
$data[19,1] = 516
$data[19,2] = 18.08
$data[19,3] = 22.88
$data[19,4] = 23.8
$data[19,5] = 73.59999999999999
$data[20,0] = "2025-06-09 09:02:17"
$data[20,1] = 516
$data[20,2] = 18.08
$data[20,3] = 20.13
$data[20,4] = 23.8
$data[20,5] = 73.7
$data[21,0] = "2025-06-09 09:02:19"
$data[21,1] = 516
$data[21,2] = 18.08
$data[21,3] = 24.04
$data[21,4] = 23.8
$data[21,5] = 73.7
$data[22,0] = "2025-06-09 09:02:20"
$data[22,1] = 516
$data[22,2] = 18.08
$data[22,3] = 789.89
$data[22,4] = 23.8
$data[22,5] = 73.7
$data[23,0] = "2025-06-09 09:02:21"
$data[23,1] = 516
$data[23,2] = 18.08
$data[23,3] = 789.84
$data[23,4] = 23.8
$data[23,5] = 73.7
$data[24,0] = "2025-06-09 09:02:23"
$data[24,1] = 516
$data[24,2] = 18.08
$data[24,3] = 17.61
$data[24,4] = 23.8
$data[24,5] = 73.8
$data[25,0] = "2025-06-09 09:02:24"
$data[25,1] = 516
$data[25,2] = 18.08
$data[25,3] = 21.71
$data[25,4] = 23.8
$data[25,5] = 73.8
$data[26,0] = "2025-06-09 09:02:25"
$data[26,1] = 516
$data[26,2] = 18.08
$data[26,3] = 789.77
$data[26,4] = 23.8
$data[26,5] = 73.8
$data[27,0] = "2025-06-09 09:02:26"
$data[27,1] = 516
$data[27,2] = 18.08
$data[27,3] = 22.88
$data[27,4] = 23.8
$data[27,5] = 73.8
$data[28,0] = "2025-06-09 09:02:28"
$data[28,1] = 516
$data[28,2] = 18.08
$data[28,3] = 789.74
$data[28,4] = 23.8
$data[28,5] = 73.90000000000001
$data[29,0] = "2025-06-09 09:02:29"
$data[29,1] = 516
$data[29,2] = 18.08
$data[29,3] = 22.49
$data[29,4] = 23.8
$data[29,5] = 73.90000000000001
$data[30,0] = "2025-06-09 09:02:30"
$data[30,1] = 516
$data[30,2] = 18.08
$data[30,3] = 21.64
$data[30,4] = 23.8
$data[30,5] = 73.90000000000001
$data[31,0] = "2025-06-09 09:02:32"
$data[31,1] = 516
$data[31,2] = 18.08
$data[31,3] = 789.77
$data[31,4] = 23.8
$data[31,5] = 73.90000000000001
$data[32,0] = "2025-06-09 09:02:33"
$data[32,1] = 516
$data[32,2] = 18.08
$data[32,3] = 789.87
$data[32,4] = 23.8
$data[32,5] = 74
$data[33,0] = "2025-06-09 09:02:34"
$data[33,1] = 516
$data[33,2] = 18.08
$data[33,3] = 789.89
$data[33,4] = 23.8
$data[33,5] = 74
$data[34,0] = "2025-06-09 09:02:36"
$data[34,1] = 516
$data[34,2] = 18.08
$data[34,3] = 22.2
$data[34,4] = 23.8
$data[34,5] = 74
$data[35,0] = "2025-06-09 09:02:37"
$data[35,1] = 516
$data[35,2] = 18.08
$data[35,3] = 23.94
$data[35,4] = 23.8
$data[35,5] = 74
$data[36,0] = "2025-06-09 09:02:38"
$data[36,1] = 516
$data[36,2] = 18.08
$data[36,3] = 789.89
$data[36,4] = 23.8
$data[36,5] = 74.09999999999999
$data[37,0] = "2025-06-09 09:02:39"
$data[37,1] = 516
$data[37,2] = 18.08
$data[37,3] = 29.39
$data[37,4] = 23.8
$data[37,5] = 74.09999999999999
$data[38,0] = "2025-06-09 09:02:41"
$data[38,1] = 516
$data[38,2] = 18.08
$data[38,3] = 789.75
$data[38,4] = 23.8
$data[38,5] = 74.09999999999999
$data[39,0] = "2025-06-09 09:02:42"
$data[39,1] = 516
$data[39,2] = 18.08
$data[39,3] = 21.39
$data[39,4] = 23.8
$data[39,5] = 74.09999999999999
$data[40,0] = "2025-06-09 09:02:43"
$data[40,1] = 516
$data[40,2] = 18.08
$data[40,3] = 789.8200000000001
$data[40,4] = 23.8
$data[40,5] = 74.09999999999999
$data[41,0] = "2025-06-09 09:02:45"
$data[41,1] = 516
$data[41,2] = 18.08
$data[41,3] = 22.53
$data[41,4] = 23.8
$data[41,5] = 74.09999999999999
$data[42,0] = "2025-06-09 09:02:46"
$data[42,1] = 516
$data[42,2] = 18.08
$data[42,3] = 22.61
$data[42,4] = 23.8
$data[42,5] = 74.09999999999999
$data[43,0] = "2025-06-09 09:02:47"
$data[43,1] = 516
$data[43,2] = 18.08
$data[43,3] = 22.78
$data[43,4] = 23.8
$data[43,5] = 74.09999999999999
$data[44,0] = "2025-06-09 09:02:48"
$data[44,1] = 516
$data[44,2] = 18.08
$data[44,3] = 24.22
$data[44,4] = 23.8
$data[44,5] = 74.09999999999999
$data[45,0] = "2025-06-09 09:02:50"
$data[45,1] = 516
$data[45,2] = 18.08
$data[45,3] = 20.35
$data[45,4] = 23.8
$data[45,5] = 74.09999999999999
$data[46,0] = "2025-06-09 09:02:51"
$data[46,1] = 516
$data[46,2] = 18.57
$data[46,3] = 20.35
$data[46,4] = 23.8
$data[46,5] = 74.09999999999999
$data[47,0] = "2025-06-09 09:02:52"
$data[47,1] = 516
$data[47,2] = 18.08
$data[47,3] = 21.4
$data[47,4] = 23.8
$data[47,5] = 74.09999999999999
$data[48,0] = "2025-06-09 09:02:54"
$data[48,1] = 516
$data[48,2] = 18.08
$data[48,3] = 19.5
$data[48,4] = 23.8
$data[48,5] = 74.09999999999999
$data[49,0] = "2025-06-09 09:02:55"
$data[49,1] = 516
$data[49,2] = 18.08
$data[49,3] = 23.17
$data[49,4] = 23.8
$data[49,5] = 74.09999999999999
$data[50,0] = "2025-06-09 10:12:58"
$data[50,1] = 521
$data[50,2] = 23.95
$data[50,3] = 267.78
$data[50,4] = -1
$data[50,5] = -1
$data[51,0] = "2025-06-09 10:12:59"
$data[51,1] = 521
$data[51,2] = 23.95
$data[51,3] = 81.65000000000001
$data[51,4] = -1
$data[51,5] = -1
$data[52,0] = "2025-06-09 10:13:00"
$data[52,1] = 521
$data[52,2] = 23.95
$data[52,3] = 57.38
$data[52,4] = -1
$data[52,5] = -1
$data[53,0] = "2025-06-09 10:13:02"
$data[53,1] = 521
$data[53,2] = 23.95
$data[53,3] = 82.77
$data[53,4] = -1
$data[53,5] = -1
$data[54,0] = "2025-06-09 10:13:03"
$data[54,1] = 521
$data[54,2] = 23.95
$data[54,3] = 267.78
$data[54,4] = -1
$data[54,5] = -1
$data[55,0] = "2025-06-09 10:13:04"
$data[55,1] = 521
$data[55,2] = 23.95
$data[55,3] = 268.07
$data[55,4] = -1
$data[55,5] = -1
$data[56,0] = "2025-06-09 10:13:06"
$data[56,1] = 521
$data[56,2] = 23.95
$data[56,3] = 266.9
$data[56,4] = -1
$data[56,5] = -1
$data[57,0] = "2025-06-09 10:13:07"
$data[57,1] = 521
$data[57,2] = 24.44
$data[57,3] = 80.48
$data[57,4] = -1
$data[57,5] = -1
$data[58,0] = "2025-06-09 10:13:08"
$data[58,1] = 521
$data[58,2] = 24.44
$data[58,3] = 266.1
$data[58,4] = -1
$data[58,5] = -1
$data[59,0] = "2025-06-09 10:13:09"
$data[59,1] = 521
$data[59,2] = 24.44
$data[59,3] = 267.67
$data[59,4] = -1
$data[59,5] = -1
$data[60,0] = "2025-06-09 10:13:11"
$data[60,1] = 521
$data[60,2] = 24.44
$data[60,3] = 267.29
$data[60,4] = -1
$data[60,5] = -1
$data[61,0] = "2025-06-09 10:13:12"
$data[61,1] = 521
$data[61,2] = 24.44
$data[61,3] = 267.31
$data[61,4] = -1
$data[61,5] = -1
$data[62,0] = "2025-06-09 10:13:13"
$data[62,1] = 521
$data[62,2] = 23.95
$data[62,3] = 267.7
$data[62,4] = -1
$data[62,5] = -1
$data[63,0] = "2025-06-09 10:13:15"
$data[63,1] = 521
$data[63,2] = 24.44
$data[63,3] = 266.97
$data[63,4] = -1
$data[63,5] = -1
$data[64,0] = "2025-06-09 10:13:16"
$data[64,1] = 521
$data[64,2] = 23.95
$data[64,3] = 87.13
$data[64,4] = -1
$data[64,5] = -1
$data[65,0] = "2025-06-09 10:13:17"
$data[65,1] = 521
$data[65,2] = 24.44
$data[65,3] = 86.79000000000001
$data[65,4] = -1
$data[65,5] = -1
$data[66,0] = "2025-06-09 10:13:18"
$data[66,1] = 521
$data[66,2] = 24.44
$data[66,3] = 269.38
$data[66,4] = -1
$data[66,5] = -1
$data[67,0] = "2025-06-09 10:13:37"
$data[67,1] = 520
$data[67,2] = 24.44
$data[67,3] = 6.14
$data[67,4] = -1
$data[67,5] = -1
$data[68,0] = "2025-06-09 10:13:52"
$data[68,1] = 518
$data[68,2] = 23.95
$data[68,3] = 278.24
$data[68,4] = -1
$data[68,5] = -1

$ws.Range("A1866:F1934").Value = $data

